$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to Text format before writing so that
# numeric-looking strings (e.g. "1.002", "319.32") are kept as literal
# text instead of being auto-converted to numbers by Excel.
$priceRange = $ws.Range("D2:D50")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "28.242.26"
$ws.Range("E2").Value = "  -2.35%  "
$ws.Range("D3").Value = "1.864.16"
$ws.Range("E3").Value = "  -1.94%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "319.32"
$ws.Range("E5").Value = "  -1.45%  "
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("D7").Value = "0.4393"
$ws.Range("E7").Value = "  -4.10%  "
$ws.Range("D8").Value = "0.3726"
$ws.Range("E8").Value = "  -2.18%  "
$ws.Range("D9").Value = "0.07556"
$ws.Range("E9").Value = "  -1.96%  "
$ws.Range("D10").Value = "0.9398"
$ws.Range("E10").Value = "  -3.47%  "
$ws.Range("D11").Value = "21.33"
$ws.Range("E11").Value = "  -2.88%  "
$ws.Range("D12").Value = "1.842.44"
$ws.Range("E12").Value = "  -3.83%  "
$ws.Range("D13").Value = "6.742"
$ws.Range("E13").Value = "  -2.46%  "
$ws.Range("D14").Value = "5.460"
$ws.Range("E14").Value = "  -3.25%  "
$ws.Range("D15").Value = "0.06886"
$ws.Range("E15").Value = "  -2.36%  "
$ws.Range("E16").Value = "  -0.19%  "
$ws.Range("D17").Value = "82.21"
$ws.Range("E17").Value = "  -1.70%  "
$ws.Range("D18").Value = "0.000009130"
$ws.Range("E18").Value = "  -3.51%  "
$ws.Range("E19").Value = "  -0.17%  "
$ws.Range("D20").Value = "16.03"
$ws.Range("E20").Value = "  -3.49%  "
$ws.Range("D21").Value = "28.223.69"
$ws.Range("E21").Value = "  -2.34%  "
$ws.Range("D22").Value = "5.170"
$ws.Range("E22").Value = "  -1.95%  "
$ws.Range("D23").Value = "10.74"
$ws.Range("E23").Value = "  -0.95%  "
$ws.Range("D24").Value = "2.085.55"
$ws.Range("E24").Value = "  -3.05%  "
$ws.Range("D25").Value = "2.035"
$ws.Range("E25").Value = "  -2.85%  "
$ws.Range("D26").Value = "154.73"
$ws.Range("E26").Value = "  -2.17%  "
$ws.Range("D27").Value = "18.48"
$ws.Range("E27").Value = "  -2.80%  "
$ws.Range("D28").Value = "5.360"
$ws.Range("E28").Value = "  -4.30%  "
$ws.Range("D29").Value = "114.40"
$ws.Range("E29").Value = "  -2.50%  "
$ws.Range("D30").Value = "1.738"
$ws.Range("E30").Value = "  -4.75%  "
$ws.Range("D31").Value = "0.09050"
$ws.Range("E31").Value = "  -2.05%  "
$ws.Range("D32").Value = "0.8026"
$ws.Range("E32").Value = "  -6.24%  "
$ws.Range("D33").Value = "4.862"
$ws.Range("E33").Value = "  -4.15%  "
$ws.Range("D34").Value = "1.169"
$ws.Range("E34").Value = "  -5.29%  "
$ws.Range("D35").Value = "2.957"
$ws.Range("E35").Value = "  -1.42%  "
$ws.Range("D36").Value = "1.002"
$ws.Range("E36").Value = "  -0.14%  "
$ws.Range("D37").Value = "1.118"
$ws.Range("E37").Value = "  -1.95%  "
$ws.Range("D38").Value = "0.05472"
$ws.Range("E38").Value = "  -3.21%  "
$ws.Range("D39").Value = "0.01954"
$ws.Range("E39").Value = "  -3.82%  "
$ws.Range("D40").Value = "2.969"
$ws.Range("E40").Value = "  +7.80%  "
$ws.Range("D41").Value = "7.153"
$ws.Range("E41").Value = "  -3.03%  "
$ws.Range("D42").Value = "0.5261"
$ws.Range("E42").Value = "  -3.83%  "
$ws.Range("D43").Value = "0.1677"
$ws.Range("E43").Value = "  -4.15%  "
$ws.Range("D44").Value = "8.777"
$ws.Range("E44").Value = "  -5.08%  "
$ws.Range("D45").Value = "2.069"
$ws.Range("E45").Value = "  +0.30%  "
$ws.Range("D46").Value = "0.06770"
$ws.Range("E46").Value = "  -0.59%  "
$ws.Range("D47").Value = "0.4886"
$ws.Range("E47").Value = "  -4.91%  "
$ws.Range("D48").Value = "0.000002536"
$ws.Range("D49").Value = "10.59"
$ws.Range("E49").Value = "  -5.45%  "
$ws.Range("D50").Value = "107.76"
$ws.Range("E50").Value = "  -2.01%  "
$ws.Range("E51").Value = "  -4.59%  "

# Restore the original (default/"Normal") cell style on column D so we
# do not leave a stray text-format style applied to those cells -
# matches the source workbook, which stores these as plain inline
# strings with no special number format.
$priceRange.Style = "Normal"
